$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clean up trailing/leading stray characters ("-" or "0") on stratum names
$ws.Range("A5").Value = "黄土"
$ws.Range("A7").Value = "黄土"
$ws.Range("A8").Value = "粉砂"
$ws.Range("A9").Value = "黏土"
$ws.Range("A11").Value = "粉砂"
$ws.Range("A12").Value = "砾石"
$ws.Range("A13").Value = "页岩"
$ws.Range("A14").Value = "粉土"
$ws.Range("A15").Value = "黏土"
$ws.Range("A16").Value = "泥岩"
$ws.Range("A17").Value = "砾石"
$ws.Range("A18").Value = "页岩"
